{"js": "// Update the worksheet date and the 20 division problems laid out in the\n// 5-column table (content appears every 4th row: rows 0, 4, 8, 12, 16).\n// Each cell/paragraph is addressed positionally (not by searching for its\n// old text) so that the new values - some of which happen to equal other\n// cells' old values - can never collide with an in-flight search.\n\n// [row, col, newText] for every populated table cell, in document order.\nconst cellUpdates = [\n  [0, 0, \"63\u00f76=10, 3\"],\n  [0, 1, \"58\u00f79=6, 4\"],\n  [0, 2, \"79\u00f79=8, 7\"],\n  [0, 3, \"75\u00f79=8, 3\"],\n  [0, 4, \"70\u00f74=17, 2\"],\n\n  [4, 0, \"15\u00f75=3, 0\"],\n  [4, 1, \"20\u00f72=10, 0\"],\n  [4, 2, \"13\u00f79=1, 4\"],\n  [4, 3, \"76\u00f72=38, 0\"],\n  [4, 4, \"70\u00f72=35, 0\"],\n\n  [8, 0, \"17\u00f73=5, 2\"],\n  [8, 1, \"59\u00f77=8, 3\"],\n  [8, 2, \"30\u00f76=5, 0\"],\n  [8, 3, \"53\u00f75=10, 3\"],\n  [8, 4, \"18\u00f72=9, 0\"],\n\n  [12, 0, \"14\u00f79=1, 5\"],\n  [12, 1, \"67\u00f78=8, 3\"],\n  [12, 2, \"37\u00f74=9, 1\"],\n  [12, 3, \"63\u00f74=15, 3\"],\n  [12, 4, \"49\u00f77=7, 0\"],\n\n  [16, 0, \"61\u00f75=12, 1\"],\n  [16, 1, \"36\u00f78=4, 4\"],\n  [16, 2, \"36\u00f72=18, 0\"],\n  [16, 3, \"12\u00f74=3, 0\"],\n  [16, 4, \"74\u00f77=10, 4\"],\n];\n\n// 1) Update the date/title paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleRange = paragraphs.items[0].getRange();\ntitleRange.insertText(\"2024-10-04 Friday\", Word.InsertLocation.replace);\n\n// 2) Update each table cell in place, preserving its run/paragraph\n// formatting by replacing text on the cell's paragraph range rather than\n// clearing/re-inserting the whole cell body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, newText] of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const cellParagraphs = cell.body.paragraphs;\n  cellParagraphs.load(\"items\");\n  await context.sync();\n\n  const cellRange = cellParagraphs.items[0].getRange();\n  cellRange.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date/title paragraph and the 20 division problems laid out in\n# the 5-column table (content appears every 4th row: rows 1, 5, 9, 13, 17 in\n# 1-based Word COM indexing). Each cell/paragraph is addressed positionally\n# (row/column, or paragraph index) rather than by searching for its old\n# text, so the new values - some of which happen to equal other cells' old\n# values - can never collide with an in-flight Find/Replace.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/title paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-04 Friday\"\n\n# 2) Update each populated table cell in place. Word COM table rows/columns\n# are 1-based, so body row R (0-based, from the OOXML) is Cell(R+1, C+1).\n$t = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @(1, 1, \"63\u00f76=10, 3\"),\n    @(1, 2, \"58\u00f79=6, 4\"),\n    @(1, 3, \"79\u00f79=8, 7\"),\n    @(1, 4, \"75\u00f79=8, 3\"),\n    @(1, 5, \"70\u00f74=17, 2\"),\n\n    @(5, 1, \"15\u00f75=3, 0\"),\n    @(5, 2, \"20\u00f72=10, 0\"),\n    @(5, 3, \"13\u00f79=1, 4\"),\n    @(5, 4, \"76\u00f72=38, 0\"),\n    @(5, 5, \"70\u00f72=35, 0\"),\n\n    @(9, 1, \"17\u00f73=5, 2\"),\n    @(9, 2, \"59\u00f77=8, 3\"),\n    @(9, 3, \"30\u00f76=5, 0\"),\n    @(9, 4, \"53\u00f75=10, 3\"),\n    @(9, 5, \"18\u00f72=9, 0\"),\n\n    @(13, 1, \"14\u00f79=1, 5\"),\n    @(13, 2, \"67\u00f78=8, 3\"),\n    @(13, 3, \"37\u00f74=9, 1\"),\n    @(13, 4, \"63\u00f74=15, 3\"),\n    @(13, 5, \"49\u00f77=7, 0\"),\n\n    @(17, 1, \"61\u00f75=12, 1\"),\n    @(17, 2, \"36\u00f78=4, 4\"),\n    @(17, 3, \"36\u00f72=18, 0\"),\n    @(17, 4, \"12\u00f74=3, 0\"),\n    @(17, 5, \"74\u00f77=10, 4\")\n)\n\nforeach ($update in $cellUpdates) {\n    $row = $update[0]\n    $col = $update[1]\n    $value = $update[2]\n    $t.Cell($row, $col).Range.Text = $value\n}\n"}
